$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the three paragraphs that need to change:
#   P4: "You can open three types of files:"
#   P5: "Soil, which will appear ..."
#   P6: "Ideal, which will show ..." + bookmark "_GoBack" + "."
# ------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p5 = $d.Paragraphs.Item(5)
$p6 = $d.Paragraphs.Item(6)

# Sanity checks so the script fails loudly if the document shape changes.
if ($p4.Range.Text -notmatch "You can open three types of files:") {
    Write-Host "Unexpected paragraph 4 text: $($p4.Range.Text)"
}
if ($p5.Range.Text -notmatch "^Soil, which will appear") {
    Write-Host "Unexpected paragraph 5 text: $($p5.Range.Text)"
}
if ($p6.Range.Text -notmatch "^Ideal, which will show") {
    Write-Host "Unexpected paragraph 6 text: $($p6.Range.Text)"
}

# 1) Rewrite paragraph 4 in place (keeps its paragraph identity/rsids).
$p4.Range.Text = "All files should contain:"

# 2) Remove the old "Soil" and "Ideal" paragraphs completely (including
#    their trailing paragraph marks and the stray bookmark they carry);
#    they will be retyped fresh further down, same as the rest of the
#    new block.
$oldRange = $d.Range($p5.Range.Start, $p6.Range.End)
$oldRange.Delete()

# ------------------------------------------------------------------
# 3) Insert the new paragraphs after paragraph 4, in order:
#      Nutrient names like "Nitrogen"
#      The unit's mg/kg, g/sqm and % (Other units do not work)
#      And values like 0.0156
#      <empty paragraph that will host the _GoBack bookmark>
#      You can open three types of files:
#      Soil, which will appear ...
#      Ideal, which will show ... in their name.
# ------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(4)

$newTexts = @(
    "Nutrient names like “Nitrogen”",
    "The unit’s mg/kg, g/sqm and % (Other units do not work)",
    "And values like 0.0156",
    "",
    "You can open three types of files:",
    "Soil, which will appear as solid grey bars at the bottom of the graph. These are what the soil is currently. These are identified by the program if they have “Soil” in their name.",
    "Ideal, which will show as green outline bars, are what the soil should be at. Max values are an indicated multiple of ideal values, these are what the soil should stay under. These are identified by the program if they have “Ideal” in their name."
)

$emptyParaForBookmark = $null

foreach ($t in $newTexts) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    if ($t -eq "") {
        # Remember this paragraph: the _GoBack bookmark goes here.
        $emptyParaForBookmark = $anchor
    }
    else {
        $anchor.Range.Text = $t
    }
}

# ------------------------------------------------------------------
# 4) Re-create the (now orphaned) "_GoBack" bookmark inside the blank
#    paragraph. A zero-width bookmark placed exactly on the paragraph
#    boundary is ambiguous, so stamp a placeholder character, bookmark
#    it, then delete the character again, leaving just the bookmark.
# ------------------------------------------------------------------
$emptyParaForBookmark.Range.Text = "X"
$start = $emptyParaForBookmark.Range.Start
$end = $start + 1
$placeholder = $d.Range($start, $end)
$d.Bookmarks.Add("_GoBack", $placeholder)
$d.Range($start, $end).Delete()

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
